$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46 (shifts existing rows 46..165 down to 47..166)
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new data record
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 45028
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 100112001
$ws.Range("G46").Value = "Berenjena"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 150
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 8000
$ws.Range("N46").Value = '$/caja 50 unidades'
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 160
$ws.Range("Q46").Value = 50
$ws.Range("R46").Value = "Hortaliza"
